# Updates cryptos list: refresh Price (col D) and Volume(1h) (col E) values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

Set-TextValue $ws.Cells.Item(2, 4) '68.810.44'
Set-TextValue $ws.Cells.Item(2, 5) '  -1.35%  '
Set-TextValue $ws.Cells.Item(3, 4) '3.498.60'
Set-TextValue $ws.Cells.Item(3, 5) '  -1.95%  '
Set-TextValue $ws.Cells.Item(4, 5) '  -0.11%  '
Set-TextValue $ws.Cells.Item(5, 4) '570.22'
Set-TextValue $ws.Cells.Item(5, 5) '  -1.53%  '
Set-TextValue $ws.Cells.Item(6, 4) '182.67'
Set-TextValue $ws.Cells.Item(6, 5) '  -3.13%  '
Set-TextValue $ws.Cells.Item(7, 5) '  -2.57%  '
Set-TextValue $ws.Cells.Item(8, 4) '3.492.73'
Set-TextValue $ws.Cells.Item(8, 5) '  -1.98%  '
Set-TextValue $ws.Cells.Item(10, 5) '  +3.56%  '
Set-TextValue $ws.Cells.Item(11, 4) '0.645'
Set-TextValue $ws.Cells.Item(11, 5) '  -2.01%  '
Set-TextValue $ws.Cells.Item(12, 4) '53.93'
Set-TextValue $ws.Cells.Item(12, 5) '  -3.17%  '
Set-TextValue $ws.Cells.Item(13, 5) '  -0.31%  '
Set-TextValue $ws.Cells.Item(14, 5) '  -1.70%  '
Set-TextValue $ws.Cells.Item(15, 4) '4.053.36'
Set-TextValue $ws.Cells.Item(15, 5) '  -2.17%  '
Set-TextValue $ws.Cells.Item(16, 4) '19.24'
Set-TextValue $ws.Cells.Item(16, 5) '  -2.30%  '
Set-TextValue $ws.Cells.Item(17, 4) '68.686.44'
Set-TextValue $ws.Cells.Item(17, 5) '  -1.42%  '
Set-TextValue $ws.Cells.Item(18, 4) '3.473.96'
Set-TextValue $ws.Cells.Item(19, 4) '12.23'
Set-TextValue $ws.Cells.Item(19, 5) '  -2.76%  '
Set-TextValue $ws.Cells.Item(20, 5) '  -1.06%  '
Set-TextValue $ws.Cells.Item(21, 4) '543.82'
Set-TextValue $ws.Cells.Item(21, 5) '  +14.76%  '
Set-TextValue $ws.Cells.Item(22, 5) '  -2.55%  '
Set-TextValue $ws.Cells.Item(23, 4) '18.95'
Set-TextValue $ws.Cells.Item(23, 5) '  -1.14%  '
Set-TextValue $ws.Cells.Item(24, 4) '4.99'
Set-TextValue $ws.Cells.Item(24, 5) '  -0.74%  '
Set-TextValue $ws.Cells.Item(25, 5) '  +0.15%  '
Set-TextValue $ws.Cells.Item(26, 4) '93.85'
Set-TextValue $ws.Cells.Item(26, 5) '  -1.66%  '
Set-TextValue $ws.Cells.Item(27, 5) '  -3.26%  '
Set-TextValue $ws.Cells.Item(28, 4) '10.76'
Set-TextValue $ws.Cells.Item(28, 5) '  -1.68%  '
Set-TextValue $ws.Cells.Item(29, 4) '9.12'
Set-TextValue $ws.Cells.Item(29, 5) '  -1.53%  '
Set-TextValue $ws.Cells.Item(30, 4) '31.53'
Set-TextValue $ws.Cells.Item(30, 5) '  -2.35%  '
Set-TextValue $ws.Cells.Item(31, 4) '7.18'
Set-TextValue $ws.Cells.Item(31, 5) '  -6.79%  '
Set-TextValue $ws.Cells.Item(32, 4) '12.53'
Set-TextValue $ws.Cells.Item(32, 5) '  +3.17%  '
Set-TextValue $ws.Cells.Item(33, 4) '64.69'
Set-TextValue $ws.Cells.Item(33, 5) '  -2.10%  '
Set-TextValue $ws.Cells.Item(34, 5) '  -4.72%  '
Set-TextValue $ws.Cells.Item(35, 4) '565.70'
Set-TextValue $ws.Cells.Item(35, 5) '  -2.86%  '
Set-TextValue $ws.Cells.Item(36, 5) '  +0.11%  '
Set-TextValue $ws.Cells.Item(37, 4) '37.71'
Set-TextValue $ws.Cells.Item(37, 5) '  -2.77%  '
Set-TextValue $ws.Cells.Item(38, 4) '0.395'
Set-TextValue $ws.Cells.Item(38, 5) '  +0.63%  '
Set-TextValue $ws.Cells.Item(39, 4) '2.95'
Set-TextValue $ws.Cells.Item(39, 5) '  +4.34%  '
Set-TextValue $ws.Cells.Item(40, 4) '0.0₃0764'
Set-TextValue $ws.Cells.Item(40, 5) '  -3.45%  '
Set-TextValue $ws.Cells.Item(41, 5) '  -3.61%  '
Set-TextValue $ws.Cells.Item(42, 5) '  -3.39%  '
Set-TextValue $ws.Cells.Item(43, 5) '  -3.87%  '
Set-TextValue $ws.Cells.Item(44, 4) '3.264.08'
Set-TextValue $ws.Cells.Item(44, 5) '  +1.43%  '
Set-TextValue $ws.Cells.Item(45, 5) '  +3.67%  '
Set-TextValue $ws.Cells.Item(46, 5) '  -3.44%  '
Set-TextValue $ws.Cells.Item(47, 4) '0.0437'
Set-TextValue $ws.Cells.Item(47, 5) '  -0.71%  '
Set-TextValue $ws.Cells.Item(48, 5) '  -2.24%  '
Set-TextValue $ws.Cells.Item(49, 4) '8.95'
Set-TextValue $ws.Cells.Item(49, 5) '  -4.90%  '
Set-TextValue $ws.Cells.Item(50, 5) '  -0.29%  '
Set-TextValue $ws.Cells.Item(51, 4) '137.54'
Set-TextValue $ws.Cells.Item(51, 5) '  +2.19%  '
